$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$newDate = Get-Date -Year 2023 -Month 9 -Day 23 -Hour 0 -Minute 0 -Second 0

for ($r = 2; $r -le 28; $r++) {
    $ws.Cells.Item($r, 3).Value = $newDate
}
